$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.465.94"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "3.148.56"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "611.32"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "143.98"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.143.73"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "5.37"
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").Value = "35.62"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "3.670.13"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("D17").Value = "64.423.69"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "3.198.08"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "6.86"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "477.05"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "14.70"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").Value = "0.725"
$ws.Range("E22").Value = "  +2.62%  "
$ws.Range("D23").Value = "7.85"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "13.70"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "84.84"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "2.81"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "8.62"
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").Value = "  +10.32%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "2.10"
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "26.63"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").Value = "2.66"
$ws.Range("E34").Value = "  -3.31%  "
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "52.66"
$ws.Range("E37").Value = "  -2.81%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0747"
$ws.Range("E38").Value = "  +4.76%  "
$ws.Range("D39").Value = "3.07"
$ws.Range("E39").Value = "  +5.87%  "
$ws.Range("D40").Value = "455.62"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").Value = "0.0397"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "8.34"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "2.859.57"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("D45").Value = "0.268"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "2.27"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "2.44"
$ws.Range("E47").Value = "  +6.15%  "
$ws.Range("D48").Value = "26.51"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "120.10"
$ws.Range("E51").Value = "  +1.38%  "
